# Bug fix on preproc omp
# Updates the "OpenMP (s)" timings (column E) on the "Pre-processing" sheet
# for the 4 datasets; dependent "Speedup" formulas in G/H recalc automatically.
# Also updates the active sheet/selection to match the author's final UI state:
#   - "Pre-processing" becomes the active (selected) tab
#   - "Main loop" selection moves to E7 (and loses its saved topLeftCell scroll)
#   - "Pre-processing" selection moves to G11

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("Main loop")
$wsPre  = $wb.Worksheets.Item("Pre-processing")

# --- Data fix: corrected OpenMP timings ---
$wsPre.Range("E4").Value = 20.86
$wsPre.Range("E5").Value = 112.02
$wsPre.Range("E6").Value = 191.4
$wsPre.Range("E7").Value = 720.07

# --- Selection / active sheet bookkeeping ---
$wsMain.Range("E7").Select()

$wsPre.Activate()
$wsPre.Range("G11").Select()
